# Add two new columns to the keytable: "STEUERBARESEINKOMMEN" (O) and
# "AMOUNT" (R). The pre-existing "HASEL"/"HASSH" columns (previously O/P)
# slide over to P/Q in the header row and in row 3, where data already
# existed in those columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
# Give the new header cells (P1, Q1, R1) the same bold style already used
# by the rest of row 1's headers (copy format from N1, which carries it).
$ws.Range("N1").Copy() | Out-Null
$ws.Range("P1:R1").PasteSpecial(-4122) | Out-Null

# Shift the existing "HASEL"/"HASSH" headers one column to the right ...
$ws.Range("P1").Value = "HASEL"
$ws.Range("Q1").Value = "HASSH"

# ... then drop in the two brand-new headers.
$ws.Range("O1").Value = "STEUERBARESEINKOMMEN"
$ws.Range("R1").Value = "AMOUNT"

# --- Row 2 ------------------------------------------------------------
# O2/P2 already held values and stay put; just append the new columns.
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -30

# --- Row 3 ------------------------------------------------------------
# O3/P3 already held values here too, but this time they get bumped one
# column to the right to make room for the new STEUERBARESEINKOMMEN value.
$oldO3 = $ws.Range("O3").Value2
$oldP3 = $ws.Range("P3").Value2
$ws.Range("P3").Value = $oldO3
$ws.Range("Q3").Value = $oldP3
$ws.Range("O3").Value = 20000000
$ws.Range("R3").Value = 50000

# --- Row 4 ------------------------------------------------------------
$ws.Range("O4").Value = 120003

# --- Row 6 ------------------------------------------------------------
$ws.Range("O6").Value = 25000
$ws.Range("R6").Value = 0

# --- Row 7 ------------------------------------------------------------
$ws.Range("O7").Value = 500000
$ws.Range("R7").Value = 3600

# --- Row 8 ------------------------------------------------------------
$ws.Range("O8").Value = 50000
$ws.Range("R8").Value = 1200

# --- Sheet view tweaks --------------------------------------------------
$ws.Range("R9").Select() | Out-Null

# --- Workbook window width ----------------------------------------------
$excel.ActiveWindow.Width = 34600
